$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Cells.Item(4, 6).Value = 1114
$wsExhibit.Cells.Item(11, 6).Value = 1433
$wsExhibit.Cells.Item(13, 6).Value = 582
$wsExhibit.Cells.Item(14, 6).Value = 1733
$wsExhibit.Cells.Item(15, 6).Value = 1789
$wsExhibit.Cells.Item(16, 6).Value = 835
$wsExhibit.Cells.Item(17, 6).Value = 266
$wsExhibit.Cells.Item(18, 6).Value = 1454
$wsExhibit.Cells.Item(19, 6).Value = 281
$wsExhibit.Cells.Item(29, 6).Value = 52

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Cells.Item(4, 6).Value = 61
$wsShow.Cells.Item(5, 6).Value = 23
$wsShow.Cells.Item(10, 6).Value = 6
$wsShow.Cells.Item(11, 6).Value = 30
$wsShow.Cells.Item(13, 6).Value = 23

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Cells.Item(6, 6).Value = 61
$wsAll.Cells.Item(8, 6).Value = 23
$wsAll.Cells.Item(12, 6).Value = 1114
$wsAll.Cells.Item(20, 6).Value = 1433
$wsAll.Cells.Item(22, 6).Value = 582
$wsAll.Cells.Item(23, 6).Value = 1733
$wsAll.Cells.Item(24, 6).Value = 1789
$wsAll.Cells.Item(25, 6).Value = 835
$wsAll.Cells.Item(26, 6).Value = 266
$wsAll.Cells.Item(27, 6).Value = 1454
$wsAll.Cells.Item(28, 6).Value = 281
$wsAll.Cells.Item(30, 6).Value = 6
$wsAll.Cells.Item(31, 6).Value = 30
$wsAll.Cells.Item(41, 6).Value = 23
$wsAll.Cells.Item(42, 6).Value = 52
